$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 164226
$ws.Range("C4").Value = 155187
$ws.Range("C5").Value = 9039
$ws.Range("C8").Value = 64.81999999999999
